# Natmi following Dr Hou advice
# Rebuild the LR-pair table for Wnt1-Fzd1 with the updated NATMI output:
# - the "ECs" sending cluster rows are added (rows 5-7)
# - the existing "FAPs" sending cluster rows (2-4) get refreshed numeric values
# - every row now also carries the updated "ECs"-derived edge weights

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs", "Wnt1", "Fzd1", "FAPs", 1, 0.3333333333333333, 0.027123, 0.081369, 0.07131444737854614, 0.07131444737854616, 3, 1, 0.668351, 2.005053, 0.04158981742241631, 0.0415898174224163, 0.018127684173, 0.163149157557, 0.002965954846054249, 0.002965954846054249),
    @("FAPs", "Wnt1", "Fzd1", "ECs", 1, 0.3333333333333333, 0.027123, 0.081369, 0.07131444737854614, 0.07131444737854616, 3, 1, 10.435983, 31.307949, 0.649405219104094, 0.649405219104094, 0.283055166909, 2.547496502181, 0.04631197432515214, 0.04631197432515215),
    @("FAPs", "Wnt1", "Fzd1", "sCs", 1, 0.3333333333333333, 0.027123, 0.081369, 0.07131444737854614, 0.07131444737854616, 3, 1, 4.965729333333333, 14.897188, 0.3090049634734898, 0.3090049634734898, 0.134685476708, 1.212169290372, 0.02203651820733976, 0.02203651820733976),
    @("ECs", "Wnt1", "Fzd1", "FAPs", 3, 1, 0.3532066666666667, 1.05962, 0.9286855526214538, 0.9286855526214538, 3, 1, 0.668351, 2.005053, 0.04158981742241631, 0.0415898174224163, 0.2360660288733333, 2.12459425986, 0.03862386257636206, 0.03862386257636205),
    @("ECs", "Wnt1", "Fzd1", "ECs", 3, 1, 0.3532066666666667, 1.05962, 0.9286855526214538, 0.9286855526214538, 3, 1, 10.435983, 31.307949, 0.649405219104094, 0.649405219104094, 3.686058768819999, 33.17452891937999, 0.6030932447789418, 0.6030932447789418),
    @("ECs", "Wnt1", "Fzd1", "sCs", 3, 1, 0.3532066666666667, 1.05962, 0.9286855526214538, 0.9286855526214538, 3, 1, 4.965729333333333, 14.897188, 0.3090049634734898, 0.3090049634734898, 1.753928705395555, 15.78535834856, 0.28696844526615, 0.28696844526615)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value = $rowVals[$c]
    }
}
